$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# row 2
$ws.Range("G2").Value = 41.0602225
$ws.Range("H2").Value = 82.120445
$ws.Range("I2").Value = 0.7047920429741388
$ws.Range("J2").Value = 0.6188016932468634
$ws.Range("M2").Value = 181.556244
$ws.Range("N2").Value = 544.668732
$ws.Range("O2").Value = 0.393453292404907
$ws.Range("P2").Value = 0.3935455037432071
$ws.Range("Q2").Value = 7454.73977490429
$ws.Range("R2").Value = 44728.43864942574
$ws.Range("S2").Value = 0.2773027497689556
$ws.Range("T2").Value = 0.2435266240859864
# row 3
$ws.Range("G3").Value = 41.0602225
$ws.Range("H3").Value = 82.120445
$ws.Range("I3").Value = 0.7047920429741388
$ws.Range("J3").Value = 0.6188016932468634
$ws.Range("M3").Value = 0.324361
$ws.Range("N3").Value = 0.648722
$ws.Range("O3").Value = 0.0007029276469155644
$ws.Range("P3").Value = 0.0004687282586276696
$ws.Range("Q3").Value = 13.3183348303225
$ws.Range("R3").Value = 53.27333932129
$ws.Range("S3").Value = 0.0004954178123326247
$ws.Range("T3").Value = 0.0002900498401114557
# row 4
$ws.Range("G4").Value = 41.0602225
$ws.Range("H4").Value = 82.120445
$ws.Range("I4").Value = 0.7047920429741388
$ws.Range("J4").Value = 0.6188016932468634
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 155.929759
$ws.Range("N4").Value = 467.789277
$ws.Range("O4").Value = 0.3379177477501335
$ws.Range("P4").Value = 0.3379969435488647
$ws.Range("Q4").Value = 6402.510598911378
$ws.Range("R4").Value = 38415.06359346827
$ws.Range("S4").Value = 0.2381617397940363
$ws.Range("T4").Value = 0.209153080980302
# row 5
$ws.Range("G5").Value = 41.0602225
$ws.Range("H5").Value = 82.120445
$ws.Range("I5").Value = 0.7047920429741388
$ws.Range("J5").Value = 0.6188016932468634
$ws.Range("M5").Value = 123.632576
$ws.Range("N5").Value = 370.897728
$ws.Range("O5").Value = 0.2679260321980438
$ws.Range("P5").Value = 0.2679888244493004
$ws.Range("Q5").Value = 5076.381078808161
$ws.Range("R5").Value = 30458.28647284896
$ws.Range("S5").Value = 0.1888321355988142
$ws.Range("T5").Value = 0.1658319383404635
# row 6
$ws.Range("I6").Value = 0.01145098916394972
$ws.Range("J6").Value = 0.01508081331502116
$ws.Range("M6").Value = 181.556244
$ws.Range("N6").Value = 544.668732
$ws.Range("O6").Value = 0.393453292404907
$ws.Range("P6").Value = 0.3935455037432071
$ws.Range("Q6").Value = 121.119619941036
$ws.Range("R6").Value = 1090.076579469324
$ws.Range("S6").Value = 0.00450542938784893
$ws.Range("T6").Value = 0.005934986272917268
# row 7
$ws.Range("I7").Value = 0.01145098916394972
$ws.Range("J7").Value = 0.01508081331502116
$ws.Range("M7").Value = 0.324361
$ws.Range("N7").Value = 0.648722
$ws.Range("O7").Value = 0.0007029276469155644
$ws.Range("P7").Value = 0.0004687282586276696
$ws.Range("Q7").Value = 0.216387385959
$ws.Range("R7").Value = 1.298324315754
$ws.Range("S7").Value = [double]"8.049216867870801E-06"
$ws.Range("T7").Value = [double]"7.068803363838842E-06"
# row 8
$ws.Range("I8").Value = 0.01145098916394972
$ws.Range("J8").Value = 0.01508081331502116
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 155.929759
$ws.Range("N8").Value = 467.789277
$ws.Range("O8").Value = 0.3379177477501335
$ws.Range("P8").Value = 0.3379969435488647
$ws.Range("Q8").Value = 104.023704894321
$ws.Range("R8").Value = 936.2133440488891
$ws.Range("S8").Value = 0.003869492467793073
$ws.Range("T8").Value = 0.005097268806708174
# row 9
$ws.Range("I9").Value = 0.01145098916394972
$ws.Range("J9").Value = 0.01508081331502116
$ws.Range("M9").Value = 123.632576
$ws.Range("N9").Value = 370.897728
$ws.Range("O9").Value = 0.2679260321980438
$ws.Range("P9").Value = 0.2679888244493004
$ws.Range("Q9").Value = 82.47764046854401
$ws.Range("R9").Value = 742.2987642168961
$ws.Range("S9").Value = 0.003068018091439843
$ws.Range("T9").Value = 0.004041489432031878
# row 10
$ws.Range("G10").Value = 7.215007666666668
$ws.Range("H10").Value = 21.645023
$ws.Range("I10").Value = 0.1238444334651152
$ws.Range("J10").Value = 0.1631016110880464
$ws.Range("M10").Value = 181.556244
$ws.Range("N10").Value = 544.668732
$ws.Range("O10").Value = 0.393453292404907
$ws.Range("P10").Value = 0.3935455037432071
$ws.Range("Q10").Value = 1309.929692391204
$ws.Range("R10").Value = 11789.36723152084
$ws.Range("S10").Value = 0.04872700009287
$ws.Range("T10").Value = 0.06418790569697388
# row 11
$ws.Range("G11").Value = 7.215007666666668
$ws.Range("H11").Value = 21.645023
$ws.Range("I11").Value = 0.1238444334651152
$ws.Range("J11").Value = 0.1631016110880464
$ws.Range("M11").Value = 0.324361
$ws.Range("N11").Value = 0.648722
$ws.Range("O11").Value = 0.0007029276469155644
$ws.Range("P11").Value = 0.0004687282586276696
$ws.Range("Q11").Value = 2.340267101767667
$ws.Range("R11").Value = 14.041602610606
$ws.Range("S11").Value = [double]"8.705367619922458E-05"
$ws.Range("T11").Value = [double]"7.645033414466739E-05"
# row 12
$ws.Range("G12").Value = 7.215007666666668
$ws.Range("H12").Value = 21.645023
$ws.Range("I12").Value = 0.1238444334651152
$ws.Range("J12").Value = 0.1631016110880464
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 155.929759
$ws.Range("N12").Value = 467.789277
$ws.Range("O12").Value = 0.3379177477501335
$ws.Range("P12").Value = 0.3379969435488647
$ws.Range("Q12").Value = 1125.034406646486
$ws.Range("R12").Value = 10125.30965981837
$ws.Range("S12").Value = 0.04184923202792298
$ws.Range("T12").Value = 0.0551278460356553
# row 13
$ws.Range("G13").Value = 7.215007666666668
$ws.Range("H13").Value = 21.645023
$ws.Range("I13").Value = 0.1238444334651152
$ws.Range("J13").Value = 0.1631016110880464
$ws.Range("M13").Value = 123.632576
$ws.Range("N13").Value = 370.897728
$ws.Range("O13").Value = 0.2679260321980438
$ws.Range("P13").Value = 0.2679888244493004
$ws.Range("Q13").Value = 892.0099836897496
$ws.Range("R13").Value = 8028.089853207745
$ws.Range("S13").Value = 0.03318114766812293
$ws.Range("T13").Value = 0.04370940902127254
# row 14
$ws.Range("G14").Value = 1.006859
$ws.Range("H14").Value = 2.013718
$ws.Range("I14").Value = 0.01728257102349843
$ws.Range("J14").Value = 0.01517395708366762
$ws.Range("M14").Value = 181.556244
$ws.Range("N14").Value = 544.668732
$ws.Range("O14").Value = 0.393453292404907
$ws.Range("P14").Value = 0.3935455037432071
$ws.Range("Q14").Value = 182.801538277596
$ws.Range("R14").Value = 1096.809229665576
$ws.Range("S14").Value = 0.0067998844704171
$ws.Range("T14").Value = 0.00597164258426978
# row 15
$ws.Range("G15").Value = 1.006859
$ws.Range("H15").Value = 2.013718
$ws.Range("I15").Value = 0.01728257102349843
$ws.Range("J15").Value = 0.01517395708366762
$ws.Range("M15").Value = 0.324361
$ws.Range("N15").Value = 0.648722
$ws.Range("O15").Value = 0.0007029276469155644
$ws.Range("P15").Value = 0.0004687282586276696
$ws.Range("Q15").Value = 0.326585792099
$ws.Range("R15").Value = 1.306343168396
$ws.Range("S15").Value = [double]"1.214839698219887E-05"
$ws.Range("T15").Value = [double]"7.112462480318516E-06"
# row 16
$ws.Range("G16").Value = 1.006859
$ws.Range("H16").Value = 2.013718
$ws.Range("I16").Value = 0.01728257102349843
$ws.Range("J16").Value = 0.01517395708366762
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 155.929759
$ws.Range("N16").Value = 467.789277
$ws.Range("O16").Value = 0.3379177477501335
$ws.Range("P16").Value = 0.3379969435488647
$ws.Range("Q16").Value = 156.999281216981
$ws.Range("R16").Value = 941.9956873018861
$ws.Range("S16").Value = 0.005840087475592308
$ws.Range("T16").Value = 0.0051287511158213
# row 17
$ws.Range("G17").Value = 1.006859
$ws.Range("H17").Value = 2.013718
$ws.Range("I17").Value = 0.01728257102349843
$ws.Range("J17").Value = 0.01517395708366762
$ws.Range("M17").Value = 123.632576
$ws.Range("N17").Value = 370.897728
$ws.Range("O17").Value = 0.2679260321980438
$ws.Range("P17").Value = 0.2679888244493004
$ws.Range("Q17").Value = 124.480571838784
$ws.Range("R17").Value = 746.883431032704
$ws.Range("S17").Value = 0.004630450680506819
$ws.Range("T17").Value = 0.00406645092109622
# row 18
$ws.Range("G18").Value = 3.486855
$ws.Range("H18").Value = 10.460565
$ws.Range("I18").Value = 0.05985129912543923
$ws.Range("J18").Value = 0.07882343226852796
$ws.Range("M18").Value = 181.556244
$ws.Range("N18").Value = 544.668732
$ws.Range("O18").Value = 0.393453292404907
$ws.Range("P18").Value = 0.3935455037432071
$ws.Range("Q18").Value = 633.0602971726199
$ws.Range("R18").Value = 5697.54267455358
$ws.Range("S18").Value = 0.023548690695615
$ws.Range("T18").Value = 0.03102060735888641
# row 19
$ws.Range("G19").Value = 3.486855
$ws.Range("H19").Value = 10.460565
$ws.Range("I19").Value = 0.05985129912543923
$ws.Range("J19").Value = 0.07882343226852796
$ws.Range("M19").Value = 0.324361
$ws.Range("N19").Value = 0.648722
$ws.Range("O19").Value = 0.0007029276469155644
$ws.Range("P19").Value = 0.0004687282586276696
$ws.Range("Q19").Value = 1.130999774655
$ws.Range("R19").Value = 6.78599864793
$ws.Range("S19").Value = [double]"4.207113285908458E-05"
$ws.Range("T19").Value = [double]"3.694677014628317E-05"
# row 20
$ws.Range("G20").Value = 3.486855
$ws.Range("H20").Value = 10.460565
$ws.Range("I20").Value = 0.05985129912543923
$ws.Range("J20").Value = 0.07882343226852796
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 155.929759
$ws.Range("N20").Value = 467.789277
$ws.Range("O20").Value = 0.3379177477501335
$ws.Range("P20").Value = 0.3379969435488647
$ws.Range("Q20").Value = 543.704459817945
$ws.Range("R20").Value = 4893.340138361505
$ws.Range("S20").Value = 0.02022481620038796
$ws.Range("T20").Value = 0.0266420791867934
# row 21
$ws.Range("G21").Value = 3.486855
$ws.Range("H21").Value = 10.460565
$ws.Range("I21").Value = 0.05985129912543923
$ws.Range("J21").Value = 0.07882343226852796
$ws.Range("M21").Value = 123.632576
$ws.Range("N21").Value = 370.897728
$ws.Range("O21").Value = 0.2679260321980438
$ws.Range("P21").Value = 0.2679888244493004
$ws.Range("Q21").Value = 431.08886578848
$ws.Range("R21").Value = 3879.79979209632
$ws.Range("S21").Value = 0.01603572109657718
$ws.Range("T21").Value = 0.02112379895270186
# row 22
$ws.Range("G22").Value = 4.822572
$ws.Range("H22").Value = 14.467716
$ws.Range("I22").Value = 0.08277866424785882
$ws.Range("J22").Value = 0.1090184929978733
$ws.Range("M22").Value = 181.556244
$ws.Range("N22").Value = 544.668732
$ws.Range("O22").Value = 0.393453292404907
$ws.Range("P22").Value = 0.3935455037432071
$ws.Range("Q22").Value = 875.568058739568
$ws.Range("R22").Value = 7880.112528656111
$ws.Range("S22").Value = 0.03256953798920042
$ws.Range("T22").Value = 0.04290373774417334
# row 23
$ws.Range("G23").Value = 4.822572
$ws.Range("H23").Value = 14.467716
$ws.Range("I23").Value = 0.08277866424785882
$ws.Range("J23").Value = 0.1090184929978733
$ws.Range("M23").Value = 0.324361
$ws.Range("N23").Value = 0.648722
$ws.Range("O23").Value = 0.0007029276469155644
$ws.Range("P23").Value = 0.0004687282586276696
$ws.Range("Q23").Value = 1.564254276492
$ws.Range("R23").Value = 9.385525658952
$ws.Range("S23").Value = [double]"5.818741167456096E-05"
$ws.Range("T23").Value = [double]"5.110004838110593E-05"
# row 24
$ws.Range("G24").Value = 4.822572
$ws.Range("H24").Value = 14.467716
$ws.Range("I24").Value = 0.08277866424785882
$ws.Range("J24").Value = 0.1090184929978733
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 155.929759
$ws.Range("N24").Value = 467.789277
$ws.Range("O24").Value = 0.3379177477501335
$ws.Range("P24").Value = 0.3379969435488647
$ws.Range("Q24").Value = 751.9824897201481
$ws.Range("R24").Value = 6767.842407481332
$ws.Range("S24").Value = 0.02797237978440095
$ws.Range("T24").Value = 0.03684791742358447
# row 25
$ws.Range("G25").Value = 4.822572
$ws.Range("H25").Value = 14.467716
$ws.Range("I25").Value = 0.08277866424785882
$ws.Range("J25").Value = 0.1090184929978733
$ws.Range("M25").Value = 123.632576
$ws.Range("N25").Value = 370.897728
$ws.Range("O25").Value = 0.2679260321980438
$ws.Range("P25").Value = 0.2679888244493004
$ws.Range("Q25").Value = 596.2269993054721
$ws.Range("R25").Value = 5366.042993749248
$ws.Range("S25").Value = 0.02217855906258288
$ws.Range("T25").Value = 0.04370940902127254
